$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 16087.571
$ws.Range("I21").Value = 5306.5
$ws.Range("J21").Value = 20400
$ws.Range("K21").Value = 5306.5
$ws.Range("L21").Value = 20400
$ws.Range("M21").Value = -4838.5
$ws.Range("N21").Value = -21336

$ws.Range("H23").Value = 16087.571
$ws.Range("I23").Value = 5306.5
$ws.Range("J23").Value = 20400
$ws.Range("K23").Value = 5306.5
$ws.Range("L23").Value = 20400
$ws.Range("M23").Value = -5072.5
$ws.Range("N23").Value = -20868

$ws.Range("H29").Value = 2333.3333
$ws.Range("J29").Value = 3250
$ws.Range("L29").Value = 9750
$ws.Range("N29").Value = -10312

$ws.Range("H38").Value = 73.5
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H58").Value = 201.66667
$ws.Range("I58").Value = 201.66667
$ws.Range("K58").Value = 605.00001
$ws.Range("M58").Value = -455.00001

$ws.Range("H132").Value = 1229.2142
$ws.Range("I132").Value = 904.44446
$ws.Range("K132").Value = 2713.33338
$ws.Range("M132").Value = -183.33338

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3985.368
$ws.Range("I32").Value = 1568.1857
$ws.Range("K32").Value = 1568.1857
$ws.Range("M32").Value = -1281.1857

$ws.Range("H61").Value = 43265.207
$ws.Range("I61").Value = 1443.3
$ws.Range("K61").Value = 1443.3
$ws.Range("M61").Value = -1231.3

$ws.Range("H136").Value = 43265.207
$ws.Range("I136").Value = 1443.3
$ws.Range("K136").Value = 4329.9
$ws.Range("M136").Value = -1779.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1159.4
$ws.Range("I5").Value = 933
$ws.Range("J5").Value = 1499
$ws.Range("K5").Value = 933
$ws.Range("L5").Value = 1499
$ws.Range("M5").Value = -820
$ws.Range("N5").Value = -1725

$ws.Range("H7").Value = 5505000.5
$ws.Range("I7").Value = 5505000.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5505000.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -5504887.5
$ws.Range("N7").ClearContents()

$ws.Range("H105").Value = 170499.83
$ws.Range("I105").Value = 501999.5
$ws.Range("K105").Value = 501999.5
$ws.Range("M105").Value = -500252.5

$ws.Range("H134").Value = 5793.8237
$ws.Range("I134").Value = 3821.0715
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 11463.2145
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -8928.2145
$ws.Range("N134").Value = -50070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2200
$ws.Range("I16").Value = 1400
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1400
$ws.Range("L16").Value = 3000
$ws.Range("N16").Value = -3574
$ws.Range("M16").Value = -1113

$ws.Range("H22").Value = 2608.4375
$ws.Range("I22").Value = 2198.8462
$ws.Range("J22").Value = 4383.3335
$ws.Range("K22").Value = 2198.8462
$ws.Range("L22").Value = 4383.3335
$ws.Range("M22").Value = -1848.8462
$ws.Range("N22").Value = -5083.3335

$ws.Range("H99").Value = 10103595
$ws.Range("I99").Value = 22224514
$ws.Range("K99").Value = 22224514
$ws.Range("M99").Value = -22223016

$ws.Range("H113").Value = 2200
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1400
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
$ws.Range("M113").Value = 770

$ws.Range("H126").Value = 10103595
$ws.Range("I126").Value = 22224514
$ws.Range("K126").Value = 66673542
$ws.Range("M126").Value = -66671072

$ws.Range("H132").Value = 1751005.4
$ws.Range("I132").Value = 2275423.5
$ws.Range("K132").Value = 6826270.5
$ws.Range("M132").Value = -6823740.5

$ws.Range("H134").Value = 2788925.8
$ws.Range("I134").Value = 3762862
$ws.Range("K134").Value = 11288586
$ws.Range("M134").Value = -11286051

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 3461571.5
$ws.Range("I9").Value = 4006002
$ws.Range("J9").Value = 2100495
$ws.Range("K9").Value = 12018006
$ws.Range("L9").Value = 6301485
$ws.Range("M9").Value = -12017782
$ws.Range("N9").Value = -6301933

$ws.Range("H132").Value = 4867.154
$ws.Range("I132").Value = 1199
$ws.Range("K132").Value = 10791
$ws.Range("M132").Value = -8261

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 16188.75
$ws.Range("J26").Value = 16188.75
$ws.Range("L26").Value = 16188.75
$ws.Range("N26").Value = -16748.75

$ws.Range("H46").Value = 20937.4
$ws.Range("J46").Value = 47523
$ws.Range("L46").Value = 47523
$ws.Range("N46").Value = -47835

$ws.Range("H50").Value = 16188.75
$ws.Range("J50").Value = 16188.75
$ws.Range("L50").Value = 16188.75
$ws.Range("N50").Value = -17184.75

$ws.Range("H70").Value = 106719.664
$ws.Range("I70").Value = 8418
$ws.Range("K70").Value = 8418
$ws.Range("M70").Value = -8148

$ws.Range("H73").Value = 106719.664
$ws.Range("I73").Value = 8418
$ws.Range("K73").Value = 8418
$ws.Range("M73").Value = -7482

$ws.Range("H113").Value = 2217.4546
$ws.Range("I113").Value = 1736.5
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 1736.5
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = 433.5
$ws.Range("N113").Value = -7840

$ws.Range("H126").Value = 4415.6
$ws.Range("I126").Value = 2399.875
$ws.Range("J126").Value = 5759.4165
$ws.Range("K126").Value = 7199.625
$ws.Range("L126").Value = 17278.2495
$ws.Range("M126").Value = -4729.625
$ws.Range("N126").Value = -22218.2495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 15664.286
$ws.Range("J22").Value = 26100
$ws.Range("L22").Value = 26100
$ws.Range("N22").Value = -26690

$ws.Range("H27").Value = 15664.286
$ws.Range("J27").Value = 26100
$ws.Range("L27").Value = 26100
$ws.Range("N27").Value = -26314

$ws.Range("H93").Value = 3397.923
$ws.Range("I93").Value = 2940
$ws.Range("J93").Value = 3481.182
$ws.Range("K93").Value = 2940
$ws.Range("L93").Value = 3481.182
$ws.Range("M93").Value = -1692
$ws.Range("N93").Value = -5977.182

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H107").Value = 7501.5
$ws.Range("J107").Value = 7501.5
$ws.Range("L107").Value = 22504.5
$ws.Range("N107").Value = -26344.5

$ws.Range("H132").Value = 2193.6
$ws.Range("I132").Value = 1783.6666
$ws.Range("K132").Value = 5350.9998
$ws.Range("M132").Value = -2820.9998
